# UndoRedoActivityDiagram.pptx - update undo/redo activity diagram text & date placeholders
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape 7 "TextBox 47": "[command commits address book]" -> "[command commits previewImageManager]"
$shp7 = $s.Shapes.Item(7)
$tr7 = $shp7.TextFrame.TextRange

$run1 = $tr7.Runs(1, 1)
$run1.Font.Size = 16

$run2 = $tr7.Runs(2, 1)
$run2.Text = "command commits previewImageManager]"
$run2.Font.Size = 14

# Resize/reposition the textbox to its final, smaller footprint (left/width stay put).
$shp7.Top = 171.85653543307086
$shp7.Height = 43.62188976377953

# --- Shape 8 "Rectangle: Rounded Corners 50": rewrite the purge-state description
$shp8 = $s.Shapes.Item(8)
$tr8 = $shp8.TextFrame.TextRange
$tr8.Delete()
$apos = [char]0x2019
$tr8.InsertAfter("Purge redundant states, then update the model" + $apos + "s previewImage by requesting for currentPreviewImageState") | Out-Null
$tr8b = $shp8.TextFrame.TextRange
$tr8b.Font.Size = 14

# --- Footer "Date Placeholder" fields on every slide layout: 6/7/2018 -> 17/10/18
$master = $p.SlideMaster
$layouts = $master.CustomLayouts
for ($l = 1; $l -le $layouts.Count; $l++) {
    $layout = $layouts.Item($l)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $dshp = $layout.Shapes.Item($i)
        if ($dshp.Name -like "Date Placeholder*") {
            $dtr = $dshp.TextFrame.TextRange
            $dtr.Text = "17/10/18"
        }
    }
}
